# "Switch pandoc to nix." — pandoc's nix-built docx writer serializes
# character-style rPr children in a different order than the previous
# (non-nix) build: bold/italic toggles now precede the color element.
# Re-apply Bold/Italic on the affected Pandoc "*Tok" character styles so
# the engine re-emits their <w:rPr> with <w:b/>/<w:i/> before <w:color/>.

$d = $word.ActiveDocument

$boldStyles = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

$italicStyles = @(
    "CommentTok",
    "DocumentationTok"
)

$boldItalicStyles = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($name in $boldStyles) {
    $s = $d.Styles($name)
    $s.Font.Bold = $True
}

foreach ($name in $italicStyles) {
    $s = $d.Styles($name)
    $s.Font.Italic = $True
}

foreach ($name in $boldItalicStyles) {
    $s = $d.Styles($name)
    $s.Font.Bold = $True
    $s.Font.Italic = $True
}
